# Generate Report for Archive
#
# The localization-status report is regenerated: every cell that used to
# read "Ready for handoff" now reads "In Translation" (Overview!E2:F2 for
# the zh-cn/de-de columns, plus the Status cell on each per-locale sheet,
# zh-cn!C2 and de-de!C2 - all four cells held the same status string).
#
# Narrowing the Status text also narrows its column - the column no
# longer needs to stay wide enough for "Ready for handoff", so it shrinks
# to fit "In Translation" on the Overview sheet (columns E and F) and on
# each locale sheet's Status column (column C).

$wb = $excel.ActiveWorkbook

$statusNew = "In Translation"
$fittedColumnWidth = 12.5   # narrowed, text-fitted width for the shorter status string

# --- Overview sheet: status shown per-locale in columns E (zh-cn) and F (de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusNew
$overview.Range("F2").Value = $statusNew
$overview.Columns.Item(5).ColumnWidth = $fittedColumnWidth
$overview.Columns.Item(6).ColumnWidth = $fittedColumnWidth

# --- zh-cn sheet: Status column (column C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusNew
$zhcn.Columns.Item(3).ColumnWidth = $fittedColumnWidth

# --- de-de sheet: Status column (column C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusNew
$dede.Columns.Item(3).ColumnWidth = $fittedColumnWidth
